# QA_517.xlsx — "Separation of qa and stg script"
#
# The "Input" worksheet's OrderId column (R2:R5) held four order IDs that
# were shared between the QA and STG scripts. This gives the QA script its
# own, separate order IDs so the two scripts no longer collide.
#
# Each new value is a purely-numeric-looking string that must stay stored as
# TEXT, matching how the existing OrderId cells are stored (General display
# format, but text content, not a number). Simply assigning .Value would let
# Excel silently coerce "51530506" into the number 51530506, so we briefly
# switch the cell to Text format while writing the value, then switch the
# format back to General to match the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

function Set-OrderId($address, $orderId) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $orderId
    $cell.NumberFormat = "General"
}

Set-OrderId "R2" "51530506"
Set-OrderId "R3" "51530507"
Set-OrderId "R4" "51530508"
Set-OrderId "R5" "51530509"
